$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shape = $s.Shapes.Item(4)
$tbl = $shape.Table

$tbl.Cell(1,1).Shape.TextFrame.TextRange.Text = "flavor"
$tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "buildType"
$tbl.Cell(1,3).Shape.TextFrame.TextRange.Text = "variant"
